$d = $word.ActiveDocument

# Locate the exact word "TEST" (case-sensitive, whole word) within the
# sentence "Le nom du fichier contenant le code pour tester les
# composants s'appelle TEST. " and replace it with the new file name.
$rng = $d.Content
$found = $rng.Find.Execute("TEST", $true, $true, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $replacement = "code_enigme2_Mathys"

    # Replace the found range's text. After this, $rng.Start/$rng.End
    # automatically re-collapse to bound exactly the freshly inserted text.
    $rng.Text = $replacement

    # Re-anchor a range over just the inserted replacement text.
    $mid = $d.Range($rng.Start, $rng.End)

    # Toggling a character property on and back off forces Word to
    # materialize the selection as its own run (split away from the
    # runs before/after it) even though the final formatting ends up
    # identical to its neighbours - matching how Word itself keeps a
    # freshly-typed replacement in a separate run from the surrounding
    # text it was typed in between.
    $mid.Bold = 1
    $mid.Bold = 0

    Write-Host "Replaced 'TEST' with '$replacement'"
} else {
    Write-Host "Target text 'TEST' not found"
}
